# Grading: fill in "Points for grading" (column E) to match column D
# (full marks) for the Generic and Customer Class sections.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Generic section (rows 3-6)
$ws.Range("E3").Value = 1
$ws.Range("E4").Value = 2
$ws.Range("E5").Value = 2
$ws.Range("E6").Value = 2

# Customer Class section (rows 10-14)
$ws.Range("E10").Value = 2
$ws.Range("E11").Value = 2
$ws.Range("E12").Value = 2
$ws.Range("E13").Value = 2
$ws.Range("E14").Value = 2

$excel.Calculate() | Out-Null

# Update selection / view to land on E15, matching the grader's last position
$ws.Activate() | Out-Null
$ws.Range("E15").Select() | Out-Null
